$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.829.89'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '3.448.21'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''579.51'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = '''147.82'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.480'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '''7.94'
$ws.Range("E9").Value = '  +2.89%  '
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("D12").Value = '4.040.43'
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("D14").Value = '''28.30'
$ws.Range("E14").Value = '  -4.81%  '
$ws.Range("D15").Value = '3.447.98'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("E16").Value = '  -1.33%  '
$ws.Range("D17").Value = '62.910.08'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").Value = '''6.45'
$ws.Range("E18").Value = '  +2.66%  '
$ws.Range("D19").Value = '''14.67'
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").Value = '''9.12'
$ws.Range("E20").Value = '  -2.53%  '
$ws.Range("D21").Value = '''388.69'
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").Value = '''0.563'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").Value = '''74.93'
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("D25").Value = '3.592.40'
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("D26").Value = '''0.0000114'
$ws.Range("E26").Value = '  -3.49%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").Value = '''7.63'
$ws.Range("E28").Value = '  -2.57%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '''8.03'
$ws.Range("E30").Value = '  -3.22%  '
$ws.Range("E31").Value = '  -1.85%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '''1.35'
$ws.Range("E33").Value = '  -6.92%  '
$ws.Range("D34").Value = '''23.29'
$ws.Range("E34").Value = '  -2.19%  '
$ws.Range("D35").Value = '''1.62'
$ws.Range("E35").Value = '  +3.22%  '
$ws.Range("D36").Value = '''5.32'
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("D37").Value = '''32.03'
$ws.Range("E37").Value = '  -1.70%  '
$ws.Range("D38").Value = '''7.01'
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("D39").Value = '''170.25'
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").Value = '3.485.96'
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("D41").Value = '''0.0783'
$ws.Range("E41").Value = '  +1.87%  '
$ws.Range("D42").Value = '''0.792'
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("E44").Value = '  -1.58%  '
$ws.Range("D45").Value = '''4.35'
$ws.Range("E45").Value = '  -3.31%  '
$ws.Range("E46").Value = '  -2.98%  '
$ws.Range("D47").Value = '2.568.94'
$ws.Range("E47").Value = '  -2.12%  '
$ws.Range("D48").Value = '''6.89'
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("D49").Value = '''2.26'
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("D50").Value = '''22.67'
$ws.Range("E50").Value = '  -4.35%  '
$ws.Range("E51").Value = '  +0.07%  '
